# Update the cryptos list on the active sheet to reflect the latest prices
# pulled from coinranking.com, as produced by the scheduled GitHub Actions
# refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-33: only the Price (D) and Volume(1h) (E) columns change ------
# Row -> @(NewPrice, NewVolume)   (NewPrice = $null means "leave unchanged")
$priceVolUpdates = @{
    2  = @('26.524.07', '  +0.09%  ')
    3  = @('1.738.32',  '  +0.19%  ')
    4  = @('0.9993',    '  -0.08%  ')
    5  = @('247.12',    '  +1.43%  ')
    6  = @('0.9995',    '  -0.10%  ')
    7  = @('0.4908',    '  +2.48%  ')
    8  = @('0.2668',    '  +0.26%  ')
    9  = @('0.06319',   '  +1.58%  ')
    10 = @('1.733.93',  '  -0.03%  ')
    11 = @('0.07039',   '  -1.10%  ')
    12 = @('15.74',     '  +0.30%  ')
    13 = @('4.615',     '  +2.02%  ')
    14 = @('0.6119',    '  +0.09%  ')
    15 = @('77.40',     '  +0.77%  ')
    16 = @('0.9996',    '  -0.08%  ')
    17 = @('0.000007426', '  +7.77%  ')
    18 = @('26.521.22', '  +0.04%  ')
    19 = @('0.9996',    '  -0.12%  ')
    20 = @($null,       '  -1.55%  ')
    21 = @('1.953.18',  '  -0.32%  ')
    22 = @('4.573',     '  +0.12%  ')
    23 = @('8.720',     '  -1.75%  ')
    24 = @('5.251',     '  -1.47%  ')
    25 = @('140.96',    '  +3.68%  ')
    26 = @('15.47',     '  +0.86%  ')
    27 = @('1.417',     '  +1.27%  ')
    28 = @($null,       '  -1.49%  ')
    29 = @('107.94',    '  +1.31%  ')
    30 = @('4.045',     '  +1.79%  ')
    31 = @('0.08065',   '  +1.59%  ')
    32 = @('3.724',     '  +0.47%  ')
    33 = @('0.04590',   '  +0.84%  ')
}

foreach ($r in $priceVolUpdates.Keys) {
    $vals = $priceVolUpdates[$r]
    $newPrice = $vals[0]
    $newVol = $vals[1]
    if ($null -ne $newPrice) {
        $ws.Cells.Item($r, 4).Value = $newPrice
    }
    $ws.Cells.Item($r, 5).Value = $newVol
}

# --- Rows 34-51: a new coin (Frax) was inserted at row 34, pushing the ----
# --- remaining coins down by one and dropping the last entry (Aave) ------
# Row -> @(Coin, Link, Price, Volume(1h))
$coinUpdates = @{
    34 = @('Frax',             'https://coinranking.com/coin/KfWtaeV1W+frax-frax',                         '0.9993',  '  -0.10%  ')
    35 = @('HuobiToken',       'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht',                  '2.607',   '  -0.23%  ')
    36 = @('ARBITRUM',         'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb',                       '1.010',   '  +2.14%  ')
    37 = @('ImmutableX',       'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx',                     '0.6375',  '  +0.41%  ')
    38 = @('TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt',               '0.8960',  '  -3.95%  ')
    39 = @('RenderToken',      'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr',                '2.019',   '  +1.96%  ')
    40 = @('MXToken',          'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx',                      '2.400',   '  +0.05%  ')
    41 = @('PaxDollar',        'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp',                     '1.005',   '  -0.01%  ')
    42 = @('VeChain',          'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet',                    '0.01506', '  -0.09%  ')
    43 = @('Quant',            'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt',                      '102.79',  '  -6.91%  ')
    44 = @('FraxShare',        'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs',                      '5.400',   '  -5.13%  ')
    45 = @('TheSandbox',       'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand',                    '0.3903',  '  +0.30%  ')
    46 = @('Aptos',            'https://coinranking.com/coin/HGYj5JCv5+aptos-apt',                          '6.901',   '  +0.02%  ')
    47 = @('Algorand',         'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo',                  '0.1188',  '  -0.18%  ')
    48 = @('Cronos',           'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro',                      '0.05396', '  +1.12%  ')
    49 = @('EnergySwap',       'https://coinranking.com/coin/SbWqqTui-+energyswap-ens',                     '7.827',   '  -1.17%  ')
    50 = @('Elrond',           'https://coinranking.com/coin/omwkOTglq+elrond-egld',                        '30.56',   '  -0.68%  ')
    51 = @('NEARProtocol',     'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near',                  '1.268',   '  +1.02%  ')
}

foreach ($r in $coinUpdates.Keys) {
    $vals = $coinUpdates[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]
    $ws.Cells.Item($r, 3).Value = $vals[1]
    $ws.Cells.Item($r, 4).Value = $vals[2]
    $ws.Cells.Item($r, 5).Value = $vals[3]
}
